$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DDF")

# Make the DDF sheet the active one (updates workbook.xml bookViews/activeTab
# and flips tabSelected from the previously-active sheet to this one).
$ws.Activate()

# Move the old row-2 detail (problem_user / abcd / password do not match)
# down to row 17, then build the new row 2 / row 3 content.
$ws.Range("A17").Value = "problem_user"
$ws.Range("B17").Value = "abcd"
$ws.Range("C17").Value = "password do not match"

$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()

$ws.Range("A2").Value = "Sauce Labs Onesie"
$ws.Range("A2").Font.Color = 1973527
$ws.Range("A2").Font.Name = "Calibri Light"
$ws.Range("A2").Font.ThemeFont = 2

$ws.Range("A3").Value = 7.99

$ws.Range("I15").Select()
